$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valentin")

# Add two new wishlist rows (24 and 25) below the existing data.
$ws.Range("A24").Value = "8bitdo Arcade Stick"
$ws.Range("B24").Value = "https://www.digitec.ch/im/productimages/2/6/7/4/8/9/8/0/2/6/7/8/8/4/3/0/0/2/7/5ca5b75d-0beb-4a12-a783-031962cefaec_cropped.jpg?impolicy=ProductTileImage&resizeWidth=648&resizeHeight=486&cropWidth=648&cropHeight=486&resizeType=downsize&quality=high"
$ws.Range("C24").Value = "https://www.digitec.ch/en/s1/product/8bitdo-arcade-stick-switch-lite-switch-pc-game-controllers-13747515"
$ws.Range("D24").Value = "78.60 CHF"

$ws.Range("A25").Value = "Dire Straits Live 1978 - 1992 LP"
$ws.Range("B25").Value = "https://superdeluxeedition.com/wp-content/uploads/2023/09/vinyl_spread-1.jpg"
$ws.Range("C25").Value = "https://www.amazon.de/dp/B0C92W29LD?tag=sdepcwde-21&linkCode=ogi&th=1&psc=1"
$ws.Range("D25").Value = "220 EUR"

# Match the cursor selection seen in the final workbook.
$ws.Range("F23").Select()
